$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2241992882562278
$ws.Range("C2").Value = 0.4911032028469751
$ws.Range("J2").Value = 0.01423487544483986
$ws.Range("O2").Value = 0.003558718861209964
$ws.Range("P2").Value = 0.1637010676156584
$ws.Range("S2").Value = 0.103202846975089
$ws.Range("B3").Value = 0.006896551724137931
$ws.Range("C3").Value = 0.02758620689655172
$ws.Range("J3").Value = 0.06206896551724138
$ws.Range("P3").Value = 0.6620689655172414
$ws.Range("S3").Value = 0.2413793103448276
$ws.Range("J4").Value = 0.06896551724137931
$ws.Range("P4").Value = 0.6206896551724138
$ws.Range("S4").Value = 0.3103448275862069
$ws.Range("B6").Value = 0.04273504273504274
$ws.Range("D6").Value = 0.004273504273504274
$ws.Range("E6").Value = 0.004273504273504274
$ws.Range("F6").Value = 0.07692307692307693
$ws.Range("J6").Value = 0.2905982905982906
$ws.Range("O6").Value = 0.02991452991452992
$ws.Range("Q6").Value = 0.1196581196581197
$ws.Range("R6").Value = 0.05128205128205128
$ws.Range("S6").Value = 0.3803418803418803
$ws.Range("B7").Value = 0.1310043668122271
$ws.Range("D7").Value = 0.01310043668122271
$ws.Range("F7").Value = 0.06550218340611354
$ws.Range("J7").Value = 0.1004366812227074
$ws.Range("O7").Value = 0.01746724890829694
$ws.Range("Q7").Value = 0.1921397379912664
$ws.Range("R7").Value = 0.08296943231441048
$ws.Range("S7").Value = 0.3973799126637554
$ws.Range("B8").Value = 0.09523809523809523
$ws.Range("D8").Value = 0.01731601731601732
$ws.Range("F8").Value = 0.08008658008658008
$ws.Range("J8").Value = 0.09956709956709957
$ws.Range("O8").Value = 0.03463203463203463
$ws.Range("Q8").Value = 0.1471861471861472
$ws.Range("R8").Value = 0.119047619047619
$ws.Range("S8").Value = 0.4069264069264069
$ws.Range("B9").Value = 0.1518987341772152
$ws.Range("D9").Value = 0.006329113924050633
$ws.Range("F9").Value = 0.03164556962025317
$ws.Range("J9").Value = 0.1139240506329114
$ws.Range("O9").Value = 0.02531645569620253
$ws.Range("Q9").Value = 0.1582278481012658
$ws.Range("R9").Value = 0.06962025316455696
$ws.Range("S9").Value = 0.4430379746835443
$ws.Range("B10").Value = 0.09722222222222222
$ws.Range("D10").Value = 0.01388888888888889
$ws.Range("E10").Value = 0.000925925925925926
$ws.Range("F10").Value = 0.07314814814814814
$ws.Range("J10").Value = 0.1212962962962963
$ws.Range("O10").Value = 0.03333333333333333
$ws.Range("Q10").Value = 0.1851851851851852
$ws.Range("R10").Value = 0.08981481481481482
$ws.Range("S10").Value = 0.3851851851851852
$ws.Range("G11").Value = 0.1475409836065574
$ws.Range("J11").Value = 0.09562841530054644
$ws.Range("K11").Value = 0.2185792349726776
$ws.Range("L11").Value = 0.5273224043715847
$ws.Range("S11").Value = 0.01092896174863388
$ws.Range("G12").Value = 0.7352941176470589
$ws.Range("J12").Value = 0.196078431372549
$ws.Range("K12").Value = 0.01470588235294118
$ws.Range("L12").Value = 0.0392156862745098
$ws.Range("S12").Value = 0.01470588235294118
$ws.Range("G13").Value = 0.6296296296296297
$ws.Range("J13").Value = 0.3148148148148148
$ws.Range("S13").Value = 0.05555555555555555
$ws.Range("F15").Value = 0.03187250996015936
$ws.Range("H15").Value = 0.1872509960159363
$ws.Range("I15").Value = 0.04382470119521913
$ws.Range("J15").Value = 0.2749003984063745
$ws.Range("K15").Value = 0.06772908366533864
$ws.Range("M15").Value = 0.0199203187250996
$ws.Range("N15").Value = 0.00398406374501992
$ws.Range("O15").Value = 0.05577689243027888
$ws.Range("S15").Value = 0.3147410358565737
$ws.Range("F16").Value = 0.01290322580645161
$ws.Range("H16").Value = 0.2
$ws.Range("I16").Value = 0.08387096774193549
$ws.Range("J16").Value = 0.3612903225806451
$ws.Range("K16").Value = 0.1161290322580645
$ws.Range("M16").Value = 0.03225806451612903
$ws.Range("O16").Value = 0.05806451612903226
$ws.Range("S16").Value = 0.1354838709677419
$ws.Range("F17").Value = 0.02785515320334262
$ws.Range("H17").Value = 0.2089136490250696
$ws.Range("I17").Value = 0.06963788300835655
$ws.Range("J17").Value = 0.3342618384401114
$ws.Range("K17").Value = 0.1587743732590529
$ws.Range("M17").Value = 0.01671309192200557
$ws.Range("O17").Value = 0.08077994428969359
$ws.Range("S17").Value = 0.1030640668523677
$ws.Range("F18").Value = 0.015625
$ws.Range("H18").Value = 0.2135416666666667
$ws.Range("I18").Value = 0.09895833333333333
$ws.Range("K18").Value = 0.08854166666666667
$ws.Range("M18").Value = 0.03125
$ws.Range("O18").Value = 0.0625
$ws.Range("S18").Value = 0.15625
$ws.Range("F19").Value = 0.02786885245901639
$ws.Range("H19").Value = 0.2213114754098361
$ws.Range("I19").Value = 0.07377049180327869
$ws.Range("J19").Value = 0.3278688524590164
$ws.Range("K19").Value = 0.1385245901639344
$ws.Range("N19").Value = 0.000819672131147541
$ws.Range("O19").Value = 0.07868852459016394
$ws.Range("S19").Value = 0.1049180327868852
